# Fruta / hortaliza, semanal
# Insert a new weekly price record as a new row 184 (shifting the existing
# rows 184-214 down to 185-215), matching the data feed's new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 184..214 down by one row, creating a blank row 184.
$ws.Rows(184).Insert()

# Populate the new row 184 with the new record.
$ws.Cells.Item(184, 1).Value2 = 5
$ws.Cells.Item(184, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(184, 3).Value2 = "Maule"
$ws.Cells.Item(184, 4).Value2 = 44491
$ws.Cells.Item(184, 5).Value2 = 7
$ws.Cells.Item(184, 6).Value2 = 100112032
$ws.Cells.Item(184, 7).Value2 = "Zapallo italiano"
$ws.Cells.Item(184, 8).Value2 = "Sin especificar"
$ws.Cells.Item(184, 9).Value2 = "Primera"
$ws.Cells.Item(184, 10).Value2 = 200
$ws.Cells.Item(184, 11).Value2 = 13000
$ws.Cells.Item(184, 12).Value2 = 13000
$ws.Cells.Item(184, 13).Value2 = 13000
$ws.Cells.Item(184, 14).Value2 = "`$/caja 60 unidades"
$ws.Cells.Item(184, 15).Value2 = "Región del Maule"
$ws.Cells.Item(184, 16).Value2 = 217
$ws.Cells.Item(184, 17).Value2 = 60
$ws.Cells.Item(184, 18).Value2 = "Hortaliza"
